# Auto-generated edit script: updates crypto price/volume data
# to reflect the commit "Updated cryptos list on Fri Aug 25 05:58:00 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.177.49'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.660.57'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.52'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5204'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2635'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06283'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.82'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.74%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.432'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.643.28'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.886.42'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5432'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.53'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.206.67'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.36'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.09'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.055'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.23%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.01'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1229'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.182'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.08'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.403'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05989'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -4.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.556'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.612'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9680'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.418'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.772'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.43%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5678'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.013'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01598'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8569'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.015.99'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -7.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.66'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.801.46'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.13%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈108'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -4.22%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.86'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.85%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.980'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05171'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.455'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.29%  '
